$d = $word.ActiveDocument

# The document ends with three paragraphs under the "Os pontos de função..."
# bullet item:
#   9.  "Os pontos de função..." (ListParagraph, numPr)
#   10. empty ListParagraph (no numPr)
#   11. empty ListParagraph (numPr)
#
# The edit removes the two trailing empty paragraphs (10 and 11), leaving the
# "Os pontos de função..." paragraph as the last paragraph in the document.

# Step 1: remove paragraph 10 (the empty ListParagraph without numPr).
# It is not the document's final paragraph, so its whole range (mark
# included) can be deleted outright.
$d.Paragraphs(10).Range.Delete()

# Step 2: the former paragraph 11 (now the last paragraph, still an empty
# ListParagraph with numPr) cannot be deleted outright because Word never
# deletes the document's final paragraph mark. Instead, merge the
# "Os pontos de função..." paragraph (now paragraph 9) into it by deleting
# just its trailing paragraph mark, which folds its text forward into the
# last paragraph.
$p = $d.Paragraphs(9)
$markRange = $d.Range($p.Range.End - 1, $p.Range.End)
$markRange.Delete()
